$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had age-bucket columns labelled with bare numbers
# (n_NA, n_2, n_3, n_4, n_5, n_6, propn_NA, propn_2, propn_3, propn_4,
# propn_5, propn_6). The final mapping iteration drops the "NA" bucket
# entirely (columns U and AA) and renames the remaining buckets to the
# "age_" naming convention (n_age_2 .. n_age_6, propn_age_2 .. propn_age_6).
#
# Deleting the higher-numbered column (AA, the old "propn_NA" column)
# first means the lower-numbered column (U, the old "n_NA" column)
# keeps its original index for the second delete, so both removals line
# up with the original layout instead of a post-shift one.
$ws.Columns.Item(27).Delete()
$ws.Columns.Item(21).Delete()

# Rename the remaining (now left-shifted) headers to the new convention.
$ws.Range("U1").Value = "n_age_2"
$ws.Range("V1").Value = "n_age_3"
$ws.Range("W1").Value = "n_age_4"
$ws.Range("X1").Value = "n_age_5"
$ws.Range("Y1").Value = "n_age_6"
$ws.Range("Z1").Value = "propn_age_2"
$ws.Range("AA1").Value = "propn_age_3"
$ws.Range("AB1").Value = "propn_age_4"
$ws.Range("AC1").Value = "propn_age_5"
$ws.Range("AD1").Value = "propn_age_6"
